$d = $word.ActiveDocument

function Add-RoleDate([string]$roleText, [string]$dateText) {
    $r = $d.Content
    $found = $r.Find.Execute($roleText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Could not find role text: $roleText"
    }
    $r.Collapse(0)
    $r.InsertAfter(" ")
    $r.Collapse(0)
    $r.InsertAfter($dateText)
}

Add-RoleDate "Head of Portfolio Architecture and Engineering" "(Mar 2023 - Present)"
Add-RoleDate "Cloud Practice Director" "(Apr 2020 - Mar 2023)"
Add-RoleDate "Technical Principal (CloudOps)" "(Aug 2019 - Apr 2020)"
Add-RoleDate "Azure Practice Lead" "(Jul 2018 - Aug 2019)"
Add-RoleDate "Windows Cloud Practice Lead - Bashton Ltd" "(Jul 2016 - Jul 2018)"
